# Applies the "calculation of new indicators" edit to the workbook.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: SCHEME_MEASURES  -> rename indicator codes MQMS0x -> MQME00x
# ---------------------------------------------------------------------------
$wsScheme = $wb.Worksheets.Item("SCHEME_MEASURES")
$wsScheme.Range("A2").Value = "MQME001"
$wsScheme.Range("A3").Value = "MQME002"
$wsScheme.Range("A4").Value = "MQME003"
$wsScheme.Range("A5").Value = "MQME004"
$wsScheme.Range("A6").Value = "MQME005"

# ---------------------------------------------------------------------------
# Sheet: METADATA_ISSUES -> renumber indicator codes in column A (rows 2-130)
# ---------------------------------------------------------------------------
$wsIssues = $wb.Worksheets.Item("METADATA_ISSUES")

# row 2 : MQME10 -> MQME012
$wsIssues.Cells.Item(2, 1).Value = "MQME012"

# rows 3-23 : MQME12 -> MQME014
for ($r = 3; $r -le 23; $r++) {
    $wsIssues.Cells.Item($r, 1).Value = "MQME014"
}

# rows 24-83 : MQME01 -> MQME008
for ($r = 24; $r -le 83; $r++) {
    $wsIssues.Cells.Item($r, 1).Value = "MQME008"
}

# rows 84-104 : MQME14 -> MQME009
for ($r = 84; $r -le 104; $r++) {
    $wsIssues.Cells.Item($r, 1).Value = "MQME009"
}

# rows 105-130 : MQME15 -> MQME010
for ($r = 105; $r -le 130; $r++) {
    $wsIssues.Cells.Item($r, 1).Value = "MQME010"
}

# ---------------------------------------------------------------------------
# Sheet: METADATA_MEASURES -> renumber/reword rows, and drop the last row
# ---------------------------------------------------------------------------
$wsMeasures = $wb.Worksheets.Item("METADATA_MEASURES")

# Wipe the old data rows (2-4) completely first so the sheet's used range
# shrinks correctly once we only repopulate rows 2-3.
$wsMeasures.Range("A2:C4").ClearContents()

$wsMeasures.Cells.Item(2, 1).Value = "MQME006"
$wsMeasures.Cells.Item(2, 2).Value = "Total number of length-required columns"
$wsMeasures.Cells.Item(2, 3).Value = 70

$wsMeasures.Cells.Item(3, 1).Value = "MQME007"
$wsMeasures.Cells.Item(3, 2).Value = "Total number of NUMBER columns"
$wsMeasures.Cells.Item(3, 3).Value = 44

# ---------------------------------------------------------------------------
# Sheet: METADATA_METRICS -> replace the 7 old rows with 11 new ones
# ---------------------------------------------------------------------------
$wsMetrics = $wb.Worksheets.Item("METADATA_METRICS")

# Wipe the old data rows (2-8) completely first, then write the full new
# set of rows (2-12) fresh - avoids any row-shift bookkeeping mistakes.
$wsMetrics.Range("A2:C8").ClearContents()

$wsMetrics.Cells.Item(2, 1).Value = "MQID001"
$wsMetrics.Cells.Item(2, 2).Value = "Table names in singular"
$wsMetrics.Cells.Item(2, 3).Value = "'95.00%"

$wsMetrics.Cells.Item(3, 1).Value = "MQID002"
$wsMetrics.Cells.Item(3, 2).Value = "Table with recommended name length"
$wsMetrics.Cells.Item(3, 3).Value = "'100.00%"

$wsMetrics.Cells.Item(4, 1).Value = "MQID003"
$wsMetrics.Cells.Item(4, 2).Value = "Columns with correct prefixes"
$wsMetrics.Cells.Item(4, 3).Value = "'84.67%"

$wsMetrics.Cells.Item(5, 1).Value = "MQID004"
$wsMetrics.Cells.Item(5, 2).Value = "Columns with recommended name size"
$wsMetrics.Cells.Item(5, 3).Value = "'100.00%"

$wsMetrics.Cells.Item(6, 1).Value = "MQID005"
$wsMetrics.Cells.Item(6, 2).Value = "Columns with comments"
$wsMetrics.Cells.Item(6, 3).Value = "'56.20%"

$wsMetrics.Cells.Item(7, 1).Value = "MQID006"
$wsMetrics.Cells.Item(7, 2).Value = "Table with standard PK prefixes"
$wsMetrics.Cells.Item(7, 3).Value = "'0.00%"

$wsMetrics.Cells.Item(8, 1).Value = "MQID007"
$wsMetrics.Cells.Item(8, 2).Value = "Table with standard FK prefixes"
$wsMetrics.Cells.Item(8, 3).Value = "'0.00%"

$wsMetrics.Cells.Item(9, 1).Value = "MQID008"
$wsMetrics.Cells.Item(9, 2).Value = "Table with standard UK prefixes"
$wsMetrics.Cells.Item(9, 3).Value = "'100.00%"

$wsMetrics.Cells.Item(10, 1).Value = "MQID009"
$wsMetrics.Cells.Item(10, 2).Value = "NUMBER columns with valid scale"
$wsMetrics.Cells.Item(10, 3).Value = "'100.00%"

$wsMetrics.Cells.Item(11, 1).Value = "MQID010"
$wsMetrics.Cells.Item(11, 2).Value = "Columns with valid num_distinct"
$wsMetrics.Cells.Item(11, 3).Value = "'100.00%"

$wsMetrics.Cells.Item(12, 1).Value = "MQID011"
$wsMetrics.Cells.Item(12, 2).Value = "Columns with valid num_nulls"
$wsMetrics.Cells.Item(12, 3).Value = "'100.00%"
